$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 and 7 (M, Density) - shift remaining rows up
$ws.Rows("6:7").Delete()

# Update C3 value (X-end, Material 2) from 6 to 2
$ws.Range("C3").Value = 2

# Update row 4 (CS-a) values to 0.5 / 0.5 with scientific number format and yellow fill
$ws.Range("B4:C4").Value = 0.5
$ws.Range("B4:C4").NumberFormat = "0.00E+00"
$ws.Range("B4:C4").Interior.Color = 65535

# Update row 5 (CS-s) values to 0.5 / 0.5 with same formatting
$ws.Range("B5:C5").Value = 0.5
$ws.Range("B5:C5").NumberFormat = "0.00E+00"
$ws.Range("B5:C5").Interior.Color = 65535

# Update selection to match the target state
$ws.Range("C6").Select()
